$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C1").Value = "UOM"
$ws.Range("C2").Value = "CU-CUBIC"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "supplier"
$ws.Range("D2").Value = "'1"
$ws.Range("E2").Value = "Supplier"

$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("E11").Select()
